$d = $word.ActiveDocument

# Highlight color used for quantitative metrics: RGB(0x2C, 0x3E, 0x50)
# Word's Font.Color / OLE_COLOR integer is packed as 0x00BBGGRR, so build it
# from the individual bytes rather than the "natural" 0xRRGGBB reading.
$metricColor = 0x2C + (0x3E * 256) + (0x50 * 65536)

function Highlight-InParagraph($paragraphIndex, $text) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $rng = $p.Range
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Font.Bold = $true
    $rng.Find.Replacement.Font.Color = $metricColor
    $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $true, $text, 2)
}

# 1) "...classification accuracy from 23% to 64%"
Highlight-InParagraph 9 "23%"
Highlight-InParagraph 9 "64%"

# 2) "Achieved 87% ... industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
Highlight-InParagraph 11 "87%"
Highlight-InParagraph 11 "71%"
Highlight-InParagraph 11 "±4.2%"
Highlight-InParagraph 11 "±2.1%"

# 3) "Wrote RFP and analyzed bids from 1,200 vendors..."
Highlight-InParagraph 31 "1,200"

# 4) "...became the $400M Polling Consortium Database ... now valued at $1B+"
Highlight-InParagraph 46 "$400M"
Highlight-InParagraph 46 "$1B"

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-InParagraph 63 "73.5%"
Highlight-InParagraph 63 "$4.7M"

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
Highlight-InParagraph 65 "87%"
Highlight-InParagraph 65 "71%"
